# ---------------------------------------------------------------
# Edit script: applies the diff changes to qcm_data.xlsx
# Sheets: users (1), history (2), questions (3)
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("users")
$ws2 = $wb.Worksheets.Item("history")
$ws3 = $wb.Worksheets.Item("questions")

# -----------------------------------------------------------
# Sheet 'users': update participant display names
# -----------------------------------------------------------
$ws1.Range("B2").Value2 = 'Meriem Ghersi'
$ws1.Range("B3").Value2 = 'mey Layadi'
$ws1.Range("B7").Value2 = 'Lina Lagab'
$ws1.Range("B12").Value2 = 'Nassim Hessas'
$ws1.Range("B13").Value2 = 'Said Mouhoun'
$ws1.Range("B14").Value2 = 'sousou'
$ws1.Range("B15").Value2 = 'meyyy'
$ws1.Range("B16").Value2 = 'merriem'
$ws1.Range("B17").Value2 = 'yara'

# -----------------------------------------------------------
# Sheet 'history': correct scores for rows 51-53 and append
# new quiz attempts as rows 54-58
# -----------------------------------------------------------
$ws2.Range("A51").Value2 = 15
$ws2.Range("B44").Copy($ws2.Range("B51"))
$ws2.Range("C51").Value2 = '3/9'

$ws2.Range("A52").Value2 = 13
$ws2.Range("B44").Copy($ws2.Range("B52"))
$ws2.Range("C52").Value2 = '2/8'

$ws2.Range("A53").Value2 = 15
$ws2.Range("B44").Copy($ws2.Range("B53"))
$ws2.Range("C53").Value2 = '5/9'

$ws2.Range("A54").Value2 = 13
$ws2.Range("B44").Copy($ws2.Range("B54"))
$ws2.Range("C54").Value2 = '2/8'

$ws2.Range("A55").Value2 = 15
$ws2.Range("B44").Copy($ws2.Range("B55"))
$ws2.Range("C55").Value2 = '2/8'

$ws2.Range("A56").Value2 = 16
$ws2.Range("B44").Copy($ws2.Range("B56"))
$ws2.Range("C56").Value2 = '2/9'

$ws2.Range("A57").Value2 = 16
$ws2.Range("B44").Copy($ws2.Range("B57"))
$ws2.Range("C57").Value2 = '2/8'

$ws2.Range("A58").Value2 = 1
$ws2.Range("B44").Copy($ws2.Range("B58"))
$ws2.Range("C58").Value2 = '2/9'

# -----------------------------------------------------------
# Sheet 'questions': fix row 10 content and fill in rows 16-28
# -----------------------------------------------------------
$ws3.Range("B10").Value2 = 'What does the map() function do in Python?'
$ws3.Range("C10").Value2 = 'Applies a function to each element of an iterable'
$ws3.Range("D10").Value2 = 'Creates a geographical map'
$ws3.Range("E10").Value2 = 'Merges two lists'
$ws3.Range("F10").Value2 = 'Applies a function to each element of an iterable'
$ws3.Range("G10").Value2 = 'The map() function takes a function and an iterable as input and applies the function to each element of the iterable.'
$ws3.Range("H10").Value2 = 'python'

$ws3.Range("B16").Value2 = 'What is the role of a Data Mart in BI?'
$ws3.Range("C16").Value2 = 'To store raw data'
$ws3.Range("D16").Value2 = 'To provide a subset of data for specific business needs'
$ws3.Range("E16").Value2 = 'To clean data'
$ws3.Range("F16").Value2 = 'To provide a subset of data for specific business needs'
$ws3.Range("G16").Value2 = 'A Data Mart is a smaller, focused version of a Data Warehouse.'
$ws3.Range("H16").Value2 = 'Business Intillegence'

$ws3.Range("B17").Value2 = 'What is the primary goal of Business Intelligence (BI)?'
$ws3.Range("C17").Value2 = 'SQL'
$ws3.Range("D17").Value2 = 'Python'
$ws3.Range("E17").Value2 = 'Power BI'
$ws3.Range("F17").Value2 = 'Power BI'
$ws3.Range("G17").Value2 = 'Power BI is a popular BI tool for data visualization and analysis.'
$ws3.Range("H17").Value2 = 'Business Intillegence'

$ws3.Range("B18").Value2 = 'What is the purpose of a KPI in BI?'
$ws3.Range("C18").Value2 = 'To measure performance'
$ws3.Range("D18").Value2 = ' To clean data'
$ws3.Range("E18").Value2 = 'To automate workflows'
$ws3.Range("F18").Value2 = 'To measure performance'
$ws3.Range("G18").Value2 = 'KPIs (Key Performance Indicators) are used to track and measure business performance.'
$ws3.Range("H18").Value2 = 'Business Intillegence'

$ws3.Range("B19").Value2 = 'Which of the following is NOT a BI tool?'
$ws3.Range("C19").Value2 = 'QlikView'
$ws3.Range("D19").Value2 = 'SAP BusinessObjects'
$ws3.Range("E19").Value2 = 'Photoshop'
$ws3.Range("F19").Value2 = 'Photoshop'
$ws3.Range("G19").Value2 = 'Photoshop is a graphic design tool, not a BI tool.'
$ws3.Range("H19").Value2 = 'Business Intillegence'

$ws3.Range("B20").Value2 = 'What is the time complexity of a binary search algorithm?'
$ws3.Range("C20").Value2 = ' O(n)'
$ws3.Range("D20").Value2 = 'O(log n)'
$ws3.Range("E20").Value2 = ' O(n^2)'
$ws3.Range("F20").Value2 = 'O(log n)'
$ws3.Range("G20").Value2 = 'Binary search divides the search space in half with each iteration, resulting in logarithmic time complexity.'
$ws3.Range("H20").Value2 = 'Algorithms'

$ws3.Range("B21").Value2 = 'Which sorting algorithm has the worst-case time complexity of O(n^2)?'
$ws3.Range("C21").Value2 = 'Merge Sort'
$ws3.Range("D21").Value2 = 'Quick Sort'
$ws3.Range("E21").Value2 = 'Bubble Sort'
$ws3.Range("F21").Value2 = 'Bubble Sort'
$ws3.Range("G21").Value2 = 'Bubble Sort has a worst-case time complexity of O(n^2).'
$ws3.Range("H21").Value2 = 'Algorithms'

$ws3.Range("A22").Value2 = 21
$ws3.Range("B22").Value2 = 'What is the purpose of a hash table?'
$ws3.Range("C22").Value2 = 'To sort data'
$ws3.Range("D22").Value2 = 'To store key-value pairs for fast lookups'
$ws3.Range("E22").Value2 = 'To perform mathematical operations'
$ws3.Range("F22").Value2 = 'To store key-value pairs for fast lookups'
$ws3.Range("G22").Value2 = 'Hash tables provide fast access to data using keys.'
$ws3.Range("H22").Value2 = 'Algorithms'

$ws3.Range("A23").Value2 = 22
$ws3.Range("B23").Value2 = 'Which algorithm is used to find the shortest path in a graph?'
$ws3.Range("C23").Value2 = 'Dijkstra''s Algorithm'
$ws3.Range("D23").Value2 = 'Bubble Sort'
$ws3.Range("E23").Value2 = 'Binary Search'
$ws3.Range("F23").Value2 = 'Dijkstra''s Algorithm'
$ws3.Range("G23").Value2 = 'Dijkstra''s Algorithm is used to find the shortest path in a weighted graph.'
$ws3.Range("H23").Value2 = 'Algorithms'

$ws3.Range("A24").Value2 = 23
$ws3.Range("B24").Value2 = 'Which data structure uses the LIFO principle?'
$ws3.Range("C24").Value2 = 'Queue'
$ws3.Range("D24").Value2 = 'Stack'
$ws3.Range("E24").Value2 = 'Linked List'
$ws3.Range("F24").Value2 = 'Stack'
$ws3.Range("G24").Value2 = 'A Stack follows the Last-In-First-Out (LIFO) principle.'
$ws3.Range("H24").Value2 = 'Algorithms'

$ws3.Range("A25").Value2 = 24
$ws3.Range("B25").Value2 = 'What is the purpose of dynamic programming?'
$ws3.Range("C25").Value2 = 'To solve problems by breaking them into smaller subproblems'
$ws3.Range("D25").Value2 = 'To sort data'
$ws3.Range("E25").Value2 = 'To perform mathematical operations'
$ws3.Range("F25").Value2 = 'To solve problems by breaking them into smaller subproblems'
$ws3.Range("G25").Value2 = 'Dynamic programming optimizes problems by storing solutions to subproblems.'
$ws3.Range("H25").Value2 = 'Algorithms'

$ws3.Range("A26").Value2 = 25
$ws3.Range("B26").Value2 = 'Which algorithm is used to detect cycles in a graph?'
$ws3.Range("C26").Value2 = 'Depth-First Search (DFS)'
$ws3.Range("D26").Value2 = 'Breadth-First Search (BFS)'
$ws3.Range("E26").Value2 = 'Binary Search'
$ws3.Range("F26").Value2 = 'Binary Search'
$ws3.Range("G26").Value2 = 'DFS can be used to detect cycles in a graph.'
$ws3.Range("H26").Value2 = 'Algorithms'

$ws3.Range("A27").Value2 = 26
$ws3.Range("B27").Value2 = 'Which algorithm is used to sort a list using the "divide and conquer" method?'
$ws3.Range("C27").Value2 = 'Bubble Sort'
$ws3.Range("D27").Value2 = 'Quick Sort'
$ws3.Range("E27").Value2 = ' Insertion Sort'
$ws3.Range("F27").Value2 = 'Quick Sort'
$ws3.Range("G27").Value2 = ' Quick Sort uses the "divide and conquer" method by dividing the list into sublists around a pivot and then sorting them recursively.'
$ws3.Range("H27").Value2 = 'Algorithms'

$ws3.Range("A28").Value2 = 27
$ws3.Range("B28").Value2 = 'What is the space complexity of the Merge Sort algorithm?'
$ws3.Range("C28").Value2 = 'O(1)'
$ws3.Range("D28").Value2 = 'O(n)'
$ws3.Range("E28").Value2 = 'O(log n)'
$ws3.Range("F28").Value2 = 'O(n)'
$ws3.Range("G28").Value2 = 'Merge Sort requires additional space to store the merged sublists, resulting in a space complexity of O(n).'
$ws3.Range("H28").Value2 = 'Algorithms'

Write-Output "Edit complete"
